# BulkUploadTemplate.xlsx - fix typo in the Award Category validation error
# message, plus the associated sample-data / active-sheet tidy-up that
# shipped in the same commit.

$wb = $excel.ActiveWorkbook

$wsData    = $wb.Worksheets.Item("Enter Your Data Here")
$wsContrib = $wb.Worksheets.Item("ReadOnly_ContributionsTypes")

# --- "Enter Your Data Here" sample row edits -------------------------------
# B1 / C1 sample value: PowerPoint -> Access
$wsData.Range("B1").Value = "Access"
$wsData.Range("C1").Value = "Access"

# A2 sample value: Code Project/Tools -> Conference (organizer)
$wsData.Range("A2").Value = "Conference (organizer)"

# New B2 sample value: Business Solutions
$wsData.Range("B2").Value = "Business Solutions"

# --- Active sheet / selection -----------------------------------------------
# The workbook used to open on "ReadOnly_ContributionsTypes" (tab 4); it
# should now open on "Enter Your Data Here" (tab 1), selection moved to C2.
$wsData.Activate()
$wsData.Range("C2").Select()

# Leaving ReadOnly_ContributionsTypes's own selection (D23) untouched;
# activating wsData above clears that sheet's tabSelected flag.
